$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 5.566218809980806
$ws.Range("H3").Value = 0.5758157389635317
$ws.Range("H4").Value = 0.7677543186180422
$ws.Range("H5").Value = 2.111324376199616
$ws.Range("B6").Value = 41
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 89
$ws.Range("E6").Value = 46
$ws.Range("F6").Value = "[458, 466, 467, 481, 553]"
$ws.Range("G6").Value = "[11, 71, 127, 190, 241, 243, 260, 318, 319, 320, 333, 334, 344, 350, 352, 357, 367, 321, 383, 386, 400, 401, 402, 403, 409, 410, 422, 423, 427, 449, 452, 454, 461, 465, 476, 478, 503, 506, 517, 520, 522]"
$ws.Range("H6").Value = 7.869481765834934
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 78
$ws.Range("E7").Value = 45
$ws.Range("F7").Value = "[25, 266, 452, 453, 457, 462, 465, 513, 514, 544]"
$ws.Range("H7").Value = 6.71785028790787
$ws.Range("H8").Value = 1.535508637236084
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 89
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "[556]"
$ws.Range("H9").Value = 1.535508637236084
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 93
$ws.Range("F10").Value = "[509]"
$ws.Range("G10").Value = "[73, 124, 141, 190, 206, 288, 350, 355, 386, 432, 433, 466, 467, 476]"
$ws.Range("H10").Value = 2.687140115163148
$ws.Range("H11").Value = 1.151631477927063
$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 91
$ws.Range("E12").Value = 23
$ws.Range("F12").Value = "[476, 499]"
$ws.Range("G12").Value = "[11, 62, 78, 121, 123, 160, 161, 190, 236, 244, 251, 259, 267, 268, 293, 350, 371, 390, 433, 442, 531]"
$ws.Range("H12").Value = 4.030710172744722
$ws.Range("H13").Value = 0.9596928982725527
$ws.Range("H14").Value = 0.5758157389635317
$ws.Range("H15").Value = 0.9596928982725527
$ws.Range("H16").Value = 1.151631477927063
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 81
$ws.Range("E17").Value = 21
$ws.Range("F17").Value = "[282, 540, 562, 569]"
$ws.Range("H17").Value = 3.262955854126679
$ws.Range("H18").Value = 1.727447216890595
$ws.Range("H19").Value = 1.919385796545105
$ws.Range("H20").Value = 4.990403071017274
$ws.Range("B21").Value = 21
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 81
$ws.Range("E21").Value = 26
$ws.Range("F21").Value = "[435, 519, 523, 556, 563]"
$ws.Range("G21").Value = "[15, 56, 91, 102, 144, 148, 160, 180, 190, 205, 269, 274, 279, 296, 350, 420, 421, 424, 436, 476, 490]"
$ws.Range("H21").Value = 4.030710172744722
$ws.Range("B22").Value = 8
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 62
$ws.Range("F22").Value = "[166, 254, 361, 371, 473]"
$ws.Range("G22").Value = "[25, 154, 190, 244, 296, 350, 476, 499]"
$ws.Range("H22").Value = 1.535508637236084
$ws.Range("H23").Value = 4.030710172744722
$ws.Range("B24").Value = 18
$ws.Range("D24").Value = 95
$ws.Range("E24").Value = 19
$ws.Range("F24").Value = "[553]"
$ws.Range("G24").Value = "[6, 107, 114, 124, 132, 170, 190, 280, 306, 350, 371, 386, 411, 451, 470, 476, 506, 536]"
$ws.Range("H24").Value = 3.45489443378119
$ws.Range("B25").Value = 48
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 75
$ws.Range("E25").Value = 64
$ws.Range("F25").Value = "[13, 34, 82, 167, 283, 321, 407, 417, 427, 429, 452, 476, 496, 511, 535, 554]"
$ws.Range("G25").Value = "[3, 4, 12, 16, 17, 18, 22, 36, 43, 64, 69, 90, 94, 97, 108, 122, 139, 152, 177, 182, 188, 190, 196, 211, 235, 244, 245, 249, 255, 290, 307, 350, 368, 372, 383, 389, 390, 398, 419, 430, 442, 443, 455, 477, 479, 517, 520, 549]"
$ws.Range("H25").Value = 9.213051823416507
$ws.Range("H26").Value = 4.414587332053743
$ws.Range("B27").Value = 116
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 95
$ws.Range("E27").Value = 122
$ws.Range("F27").Value = "[510, 547, 556, 557, 574, 575]"
$ws.Range("G27").Value = "[11, 25, 26, 28, 29, 30, 32, 33, 41, 46, 47, 48, 49, 61, 65, 71, 72, 84, 86, 92, 98, 99, 100, 118, 119, 125, 128, 138, 139, 150, 152, 154, 160, 163, 164, 165, 166, 168, 176, 177, 180, 185, 186, 187, 189, 190, 194, 195, 196, 204, 205, 217, 228, 231, 238, 239, 240, 254, 258, 261, 266, 270, 271, 272, 279, 301, 317, 328, 335, 346, 347, 348, 350, 361, 364, 365, 371, 372, 378, 379, 380, 407, 408, 412, 414, 415, 416, 424, 429, 430, 431, 434, 441, 444, 445, 446, 462, 463, 464, 475, 476, 482, 483, 484, 508, 513, 514, 515, 524, 525, 532, 533, 534, 540, 544, 545]"
$ws.Range("H27").Value = 22.26487523992322
$ws.Range("H28").Value = 0.9596928982725527

# Restore the bold header font explicitly so Excel re-derives the font family id
$ws.Range("A1:H1").Font.Name = "Calibri"

# Adjust column widths to match the refreshed report layout
$ws.Columns.Item(1).ColumnWidth = 20.67
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 13.83
$ws.Columns.Item(4).ColumnWidth = 16.33
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 62.67

# Move the active selection to the updated row of interest
$ws.Range("A17:H17").Select()
